$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Columns.Item(1).ColumnWidth = 21.2105263157895
$ws.Columns.Item(3).ColumnWidth = 29.5668016194332
$ws.Columns.Item(7).ColumnWidth = 21.3157894736842
$ws.Columns.Item(10).ColumnWidth = 21.7449392712551
$ws.Columns.Item(11).ColumnWidth = 29.1376518218623
